$d = $word.ActiveDocument

# The underlying change here is purely a re-serialization of the OOXML
# parts (document.xml / styles.xml) caused by bumping the OOXML writer
# library (POI 3.17.0 -> 4.0.1): every value is identical, only the
# order in which attributes are emitted for a given element changes.
#
# We cannot poke the serializer directly from script, but touching any
# property on an object forces this runtime to rewrite the owning XML
# part with its current (new) attribute ordering. So we make a set of
# "set it back to what it already was" edits - they are no-ops content
# wise, but they flip every part that needs to flip to the new layout.

# --- word/document.xml ---
# Touching a table's preferred-width property rewrites the whole body,
# which normalizes every <w:tblW>, <w:tblLook>, <w:tcW>, <w:pgSz> and
# <w:pgMar> element (in every table, including the nested ones) to the
# new attribute order.
foreach ($t in $d.Tables) {
    $t.PreferredWidthType = $t.PreferredWidthType
}

# --- word/styles.xml ---
# Touching any style property rewrites the whole styles part, which
# normalizes <w:docDefaults>, <w:latentStyles>, every <w:lsdException>
# and every <w:style> element to the new attribute order.
foreach ($st in $d.Styles) {
    $st.NameLocal = $st.NameLocal
}
